$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Give the new G1 header cell the same format as the other header cells
# (copy from B1, which already carries the bold/border/centred header style),
# then fix up the text of E1/F1/G1:
#   E1: "E-mail" -> "login"  (new column)
#   F1: stays "Пароль"
#   G1: new column carrying the old "E-mail" text
$ws.Range("B1").Copy()
$ws.Range("G1").PasteSpecial(-4122)

$ws.Range("E1").Value = "login"
$ws.Range("F1").Value = "Пароль"
$ws.Range("G1").Value = "E-mail"

# --- Data rows ---
$data = @(
    @(0,  "dsa",    "sda",    "asdas", "asd",    "asd",   "aaa@mail.ru"),
    @(1,  "das",    "asd",    "asd",   "asda",   "asd",   "asd"),
    @(2,  "das",    "asd",    "asd",   "asda",   "asd",   "asd"),
    @(3,  "das",    "asd",    "asd",   "asda",   "asd",   "asd"),
    @(4,  "das",    "asd",    "asd",   "asda",   "asd",   "asd"),
    @(5,  "das",    "asd",    "asd",   "asda",   "asd",   "asd"),
    @(6,  "das",    "asd",    "asd",   "asda",   "asd",   "asd"),
    @(7,  "dimas",  "dimas",  "dimas", "dimas",  "dimas", "dimas"),
    @(8,  "dimass", "dimass", "ds",    "dimass", "dimas", "dimasik"),
    @(9,  "dasd",   "sad",    "",      "",       "",      ""),
    @(10, "dimon",  "dimon",  "dimon", "dimon",  "dimon", "dimon"),
    @(11, "alex",   "alex",   "alex",  "alex",   "alex",  "alex")
)

$row = 2
foreach ($entry in $data) {
    if ($row -gt 3) {
        # New rows: clone A2's number style (bordered/bold/centred) into the
        # new A-column cell before writing its value.
        $ws.Range("A2").Copy()
        $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    }
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $ws.Cells.Item($row, 7).Value = $entry[6]
    $row = $row + 1
}
